$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-format on cells whose new values would otherwise be
# auto-coerced to numeric (losing exact decimal text / trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated Price (D) / Volume(1h) (E) text values.
$ws.Range("D2").Value = '60.907.37'
$ws.Range("E2").Value = '  -3.76%  '
$ws.Range("D3").Value = '2.918.17'
$ws.Range("E3").Value = '  -4.15%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '586.19'
$ws.Range("E5").Value = '  -1.67%  '
$ws.Range("D6").Value = '145.33'
$ws.Range("E6").Value = '  -6.30%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -2.31%  '
$ws.Range("D9").Value = '2.917.55'
$ws.Range("E9").Value = '  -4.19%  '
$ws.Range("E10").Value = '  -0.77%  '
$ws.Range("E11").Value = '  -5.08%  '
$ws.Range("E12").Value = '  -4.00%  '
$ws.Range("E13").Value = '  -4.21%  '
$ws.Range("D14").Value = '33.59'
$ws.Range("E14").Value = '  -6.34%  '
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").Value = '3.400.49'
$ws.Range("E16").Value = '  -4.18%  '
$ws.Range("D17").Value = '60.850.97'
$ws.Range("E17").Value = '  -3.78%  '
$ws.Range("D18").Value = '6.78'
$ws.Range("E18").Value = '  -4.44%  '
$ws.Range("D19").Value = '2.918.44'
$ws.Range("E19").Value = '  -4.22%  '
$ws.Range("D20").Value = '429.14'
$ws.Range("E20").Value = '  -5.84%  '
$ws.Range("E21").Value = '  -5.16%  '
$ws.Range("E22").Value = '  -2.51%  '
$ws.Range("D23").Value = '7.13'
$ws.Range("E23").Value = '  -5.54%  '
$ws.Range("D24").Value = '80.60'
$ws.Range("E24").Value = '  -3.09%  '
$ws.Range("E25").Value = '  -3.46%  '
$ws.Range("D26").Value = '10.73'
$ws.Range("E26").Value = '  -5.13%  '
$ws.Range("E27").Value = '  -3.55%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").Value = '7.18'
$ws.Range("E30").Value = '  -4.49%  '
$ws.Range("E32").Value = '  -3.66%  '
$ws.Range("D33").Value = '26.63'
$ws.Range("E33").Value = '  -3.87%  '
$ws.Range("D34").Value = '0.107'
$ws.Range("E34").Value = '  -3.88%  '
$ws.Range("D35").Value = '0.0₃0871'
$ws.Range("E35").Value = '  +0.54%  '
$ws.Range("E36").Value = '  -3.01%  '
$ws.Range("E37").Value = '  -5.21%  '
$ws.Range("D38").Value = '3.02'
$ws.Range("E38").Value = '  -6.32%  '
$ws.Range("E39").Value = '  -3.70%  '
$ws.Range("D40").Value = '49.61'
$ws.Range("E40").Value = '  -1.65%  '
$ws.Range("E41").Value = '  -6.14%  '
$ws.Range("D42").Value = '8.66'
$ws.Range("E42").Value = '  -5.23%  '
$ws.Range("E43").Value = '  -2.82%  '
$ws.Range("D44").Value = '41.02'
$ws.Range("E44").Value = '  -6.00%  '
$ws.Range("D45").Value = '378.30'
$ws.Range("E45").Value = '  -4.89%  '
$ws.Range("D46").Value = '0.0351'
$ws.Range("E46").Value = '  -3.46%  '
$ws.Range("D47").Value = '2.695.90'
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("D48").Value = '132.77'
$ws.Range("E48").Value = '  +0.28%  '
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").Value = '24.59'
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("E51").Value = '  -2.39%  '

Write-Output "Updated $(($wb.ActiveSheet.UsedRange.Rows.Count)) rows of crypto data."
